# "Cambio de link a local"
# Update the Roles/Areas table:
#  - "Miguel" -> "Miguel Carranza" (row 2, unchanged Rol/Area/Actividades)
#  - add a new row 3 for "Selena Lopez" / Soporte / Piso Productivo / Validar conteos
#  - widen the Area/Actividades columns (C:D) to fit the new text

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Write the new "Piso Productivo" area value first so it lands earlier in the
# shared-string table (matches the target ordering), then fill in row 2/3.
$ws.Range("C3").Value = "Piso Productivo"
$ws.Range("A2").Value = "Miguel Carranza"
$ws.Range("A3").Value = "Selena Lopez"
$ws.Range("B3").Value = "Soporte"
$ws.Range("D3").Value = "Validar conteos"

# Widen columns C:D to fit the longer text (bestFit-style autosize).
$ws.Range("C1").ColumnWidth = 13.3
$ws.Range("D1").ColumnWidth = 13.3

# Move the active selection down to the next empty row, like Excel does
# after entering data in the row above.
$null = $ws.Range("A4").Select()
